$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-08 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-09 Saturday", 2)

$d.Content.Find.Execute("83×48=", $true, $false, $false, $false, $false, $true, 1, $false, "64×11=", 2)
$d.Content.Find.Execute("97×53=", $true, $false, $false, $false, $false, $true, 1, $false, "63×23=", 2)
$d.Content.Find.Execute("73×51=", $true, $false, $false, $false, $false, $true, 1, $false, "38×94=", 2)
$d.Content.Find.Execute("78×89=", $true, $false, $false, $false, $false, $true, 1, $false, "55×71=", 2)
$d.Content.Find.Execute("91×35=", $true, $false, $false, $false, $false, $true, 1, $false, "91×89=", 2)

$d.Content.Find.Execute("18×11=", $true, $false, $false, $false, $false, $true, 1, $false, "15×26=", 2)
$d.Content.Find.Execute("96×88=", $true, $false, $false, $false, $false, $true, 1, $false, "70×68=", 2)
$d.Content.Find.Execute("61×71=", $true, $false, $false, $false, $false, $true, 1, $false, "65×65=", 2)
$d.Content.Find.Execute("57×84=", $true, $false, $false, $false, $false, $true, 1, $false, "63×46=", 2)
$d.Content.Find.Execute("47×34=", $true, $false, $false, $false, $false, $true, 1, $false, "83×74=", 2)

$d.Content.Find.Execute("39×58=", $true, $false, $false, $false, $false, $true, 1, $false, "80×74=", 2)
$d.Content.Find.Execute("81×40=", $true, $false, $false, $false, $false, $true, 1, $false, "84×42=", 2)
$d.Content.Find.Execute("59×25=", $true, $false, $false, $false, $false, $true, 1, $false, "31×65=", 2)
$d.Content.Find.Execute("28×79=", $true, $false, $false, $false, $false, $true, 1, $false, "30×18=", 2)
$d.Content.Find.Execute("32×85=", $true, $false, $false, $false, $false, $true, 1, $false, "67×36=", 2)

$d.Content.Find.Execute("80×96=", $true, $false, $false, $false, $false, $true, 1, $false, "40×53=", 2)
$d.Content.Find.Execute("29×32=", $true, $false, $false, $false, $false, $true, 1, $false, "86×26=", 2)
$d.Content.Find.Execute("88×38=", $true, $false, $false, $false, $false, $true, 1, $false, "89×82=", 2)
$d.Content.Find.Execute("89×93=", $true, $false, $false, $false, $false, $true, 1, $false, "73×24=", 2)
$d.Content.Find.Execute("57×30=", $true, $false, $false, $false, $false, $true, 1, $false, "41×60=", 2)

$d.Content.Find.Execute("65×99=", $true, $false, $false, $false, $false, $true, 1, $false, "54×92=", 2)
$d.Content.Find.Execute("91×95=", $true, $false, $false, $false, $false, $true, 1, $false, "60×82=", 2)
$d.Content.Find.Execute("40×15=", $true, $false, $false, $false, $false, $true, 1, $false, "90×54=", 2)
$d.Content.Find.Execute("97×94=", $true, $false, $false, $false, $false, $true, 1, $false, "32×86=", 2)
$d.Content.Find.Execute("75×67=", $true, $false, $false, $false, $false, $true, 1, $false, "86×16=", 2)
